$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The second table's row 2 has 3 cells:
#   (2,1) "{{ tema.titulo }}"
#   (2,2) <empty, italic+underline formatted>   <-- type "{{ tema.items }}" here
#   (2,3) <empty, bold formatted>, currently holds the "_GoBack" bookmark
#
# We need to:
#   1. Put the text "{{ tema.items }}" into cell (2,2), with the "_GoBack"
#      bookmark now sitting inside it (right after "{{ ", before "tema.items"),
#      exactly like Word repositions _GoBack to the most-recently-typed spot.
#   2. Remove the "_GoBack" bookmark from its old spot in cell (2,3), leaving
#      that paragraph empty as before.
# ---------------------------------------------------------------------------

$table = $d.Tables.Item(2)

# --- Step 1: cell (2,2) -----------------------------------------------------
$targetCell = $table.Cell(2, 2)
$targetRange = $targetCell.Range

$newParagraphXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="000E76DD" w:rsidRPr="00B35ACB" w:rsidRDefault="000E76DD" w:rsidP="00630DA1"><w:pPr><w:widowControl w:val="0"/><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:i/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:i/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:i/><w:u w:val="single"/></w:rPr><w:t>tema.items</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:i/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$targetRange.InsertXML($newParagraphXml)

# InsertXML drops the new paragraph in front of the (now stray) original
# empty paragraph; delete that leftover empty paragraph, leaving only the
# freshly-inserted one behind.
$table = $d.Tables.Item(2)
$targetCell = $table.Cell(2, 2)
$staleParagraph = $targetCell.Range.Paragraphs.First
$staleParagraph.Range.Delete()

# --- Step 2: cell (2,3) - drop the bookmark that used to live here ---------
$table = $d.Tables.Item(2)
$bookmarkCell = $table.Cell(2, 3)
$bookmarkRange = $bookmarkCell.Range

$cleanParagraphXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="000E76DD" w:rsidRPr="00370A7E" w:rsidRDefault="000E76DD" w:rsidP="00630DA1"><w:pPr><w:widowControl w:val="0"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:b/></w:rPr></w:pPr></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$bookmarkRange.InsertXML($cleanParagraphXml)

$table = $d.Tables.Item(2)
$bookmarkCell = $table.Cell(2, 3)
$staleBookmarkParagraph = $bookmarkCell.Range.Paragraphs.First
$staleBookmarkParagraph.Range.Delete()
